$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.614.97"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "2.619.29"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.82%  "
$ws.Range("E10").Value = "  -2.24%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "3.087.36"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "58.563.63"
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "2.632.14"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "334.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("D28").Value = "0.0₃0735"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.80%  "
$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.847"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.817"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "281.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.595"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0941"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0528"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0225"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "1.944.41"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.34%  "
